# Add 2022-Q1 data
#
# The workbook currently has three sheets: "2021-Q3", "2021-Q4", "总计"
# (总计 = "Total"/summary sheet). This change:
#   1. Turns the existing "总计" sheet into the new "2022-Q1" holdings
#      sheet (same role the 2021-Q3/2021-Q4 sheets play), populated with
#      the fund-holding data for that quarter.
#   2. Adds a brand new "总计" summary sheet (placed after "2022-Q1") that
#      lists the per-quarter holding counts/market values, now including
#      the 2022-Q1 row.

$wb = $excel.ActiveWorkbook

$oldTotal = $wb.Worksheets.Item("总计")

# --- Make a copy of the summary sheet first (while it still has the
#     original 2021-Q4 / 2021-Q3 rows + correct styles) - this copy will
#     become the new "总计" sheet. It is placed right after the original.
$oldTotal.Copy($null, $oldTotal)
$newTotal = $wb.Worksheets.Item("总计 (2)")

# Free up the "总计" name by renaming the original sheet first, then
# rename the copy into place.
$oldTotal.Name = "2022-Q1"
$newTotal.Name = "总计"

# A scratch cell (always outside any real data range) used purely to
# "donate" a plain/default cell style via PasteSpecial -> Formats. This
# lets us store numeric-looking values ("2.42", "000586", ...) as plain
# text without leaving Excel's "number stored as text" quote-prefix
# marker (and its extra style) on the cell.
$xlPasteFormats = -4122

# ======================================================================
# 1) "2022-Q1" sheet: replace the old summary content with the fund
#    holdings table (same shape as the "2021-Q3"/"2021-Q4" sheets).
# ======================================================================
$q1 = $wb.Worksheets.Item("2022-Q1")
$template = $wb.Worksheets.Item("2021-Q4")

$q1.Cells.Clear()
$template.Range("A1:H4").Copy($q1.Range("A1:H4"))
$q1.Range("A1").ClearContents()

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "'000586"
$q1.Range("C2").Value = "景顺长城中小板创业板精选股票"
$q1.Range("D2").Value = "'2.42"
$q1.Range("E2").Value = "'94.15"
$q1.Range("F2").Value = "'6.50"
$q1.Range("G2").Value = "'0.1573"
$q1.Range("H2").Value = 5

$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "'010706"
$q1.Range("C3").Value = "景顺长城景骊成长混合型证券投资基金"
$q1.Range("D3").Value = "'1.13"
$q1.Range("E3").Value = "'93.50"
$q1.Range("F3").Value = "'6.29"
$q1.Range("G3").Value = "'0.0711"
$q1.Range("H3").Value = 3

$q1.Range("A4").Value = 2
$q1.Range("B4").Value = "'260115"
$q1.Range("C4").Value = "景顺长城中小盘混合"
$q1.Range("D4").Value = "'0.96"
$q1.Range("E4").Value = "'94.00"
$q1.Range("F4").Value = "'5.79"
$q1.Range("G4").Value = "'0.0556"
$q1.Range("H4").Value = 4

# Strip the "number stored as text" quote-prefix styling that typing a
# leading apostrophe leaves behind, so the text-holding cells end up
# with a plain/default cell style (matching the rest of the workbook).
$q1.Range("Z1").Copy()
$q1.Range("B2:B4").PasteSpecial($xlPasteFormats)
$q1.Range("D2:G4").PasteSpecial($xlPasteFormats)
$q1.Range("Z1").ClearContents()

# ======================================================================
# 2) "总计" sheet: insert a new row for 2022-Q1 above the 2021-Q4 row.
# ======================================================================
$total = $wb.Worksheets.Item("总计")

$total.Rows(2).Insert()
$total.Range("B2:D2").ClearFormats()
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial($xlPasteFormats)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.28

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2

Write-Host "Added 2022-Q1 sheet and updated 总计 summary"
